# Add two new test-case worksheets ("ApproveInvoice" and "TransferInvoice")
# to the AR automation workbook, inserted right after "FirmAllLines" and
# before "Allocate", following the same layout/style as the existing
# single-row API-mode test sheets (e.g. "Allocate").

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ApproveInvoice ------------------------------------------------
$firmAllLines = $wb.Worksheets.Item("FirmAllLines")
$wsApprove = $wb.Worksheets.Add($null, $firmAllLines)
$wsApprove.Name = "ApproveInvoice"

# Re-fetch the template sheet each time (inserting a sheet shifts positions,
# so stale references can end up pointing at the wrong worksheet).
$template1 = $wb.Worksheets.Item("Allocate")
$template1.Range("A1:C2").Copy()
$wsApprove.Range("A1").PasteSpecial(-4122)   # xlPasteFormats: copy cell styles only

$wsApprove.Range("A1").Value = "Transaction Type"
$wsApprove.Range("B1").Value = "Background Processing"
$wsApprove.Range("C1").Value = "Process All Lines"
$wsApprove.Range("A2").Value = "Approve Sales Invoice"
$wsApprove.Range("B2").Value = $true
$wsApprove.Range("C2").Value = $true
$wsApprove.Range("G7").Select()

# ---- Sheet 2: TransferInvoice ----------------------------------------------
$wsTransfer = $wb.Worksheets.Add($null, $wsApprove)
$wsTransfer.Name = "TransferInvoice"

$template2 = $wb.Worksheets.Item("Allocate")
$template2.Range("A1:C2").Copy()
$wsTransfer.Range("A1").PasteSpecial(-4122)  # xlPasteFormats: copy cell styles only

$wsTransfer.Range("A1").Value = "Transaction Type"
$wsTransfer.Range("B1").Value = "Background Processing"
$wsTransfer.Range("C1").Value = "Process All Lines"
$wsTransfer.Range("A2").Value = "Transfer Sales Invoice"
$wsTransfer.Range("B2").Value = $true
$wsTransfer.Range("C2").Value = $true
$wsTransfer.Range("D8").Select()

# TransferInvoice ends up the active/selected tab, matching the authored file.
$wsTransfer.Activate()
